$d = $word.ActiveDocument

$d.Content.Find.Execute("92×52=", $true, $false, $false, $false, $false, $true, 1, $false, "85×68=", 2) | Out-Null
$d.Content.Find.Execute("11×46=", $true, $false, $false, $false, $false, $true, 1, $false, "36×50=", 2) | Out-Null
$d.Content.Find.Execute("64×69=", $true, $false, $false, $false, $false, $true, 1, $false, "64×20=", 2) | Out-Null
$d.Content.Find.Execute("47×26=", $true, $false, $false, $false, $false, $true, 1, $false, "14×79=", 2) | Out-Null
$d.Content.Find.Execute("54×87=", $true, $false, $false, $false, $false, $true, 1, $false, "44×79=", 2) | Out-Null
$d.Content.Find.Execute("15×14=", $true, $false, $false, $false, $false, $true, 1, $false, "88×82=", 2) | Out-Null
$d.Content.Find.Execute("48×12=", $true, $false, $false, $false, $false, $true, 1, $false, "71×69=", 2) | Out-Null
$d.Content.Find.Execute("56×83=", $true, $false, $false, $false, $false, $true, 1, $false, "85×28=", 2) | Out-Null
$d.Content.Find.Execute("66×19=", $true, $false, $false, $false, $false, $true, 1, $false, "96×72=", 2) | Out-Null
$d.Content.Find.Execute("90×13=", $true, $false, $false, $false, $false, $true, 1, $false, "79×95=", 2) | Out-Null
$d.Content.Find.Execute("92×25=", $true, $false, $false, $false, $false, $true, 1, $false, "47×11=", 2) | Out-Null
$d.Content.Find.Execute("95×89=", $true, $false, $false, $false, $false, $true, 1, $false, "52×93=", 2) | Out-Null
$d.Content.Find.Execute("19×35=", $true, $false, $false, $false, $false, $true, 1, $false, "80×89=", 2) | Out-Null
$d.Content.Find.Execute("25×84=", $true, $false, $false, $false, $false, $true, 1, $false, "76×56=", 2) | Out-Null
$d.Content.Find.Execute("82×23=", $true, $false, $false, $false, $false, $true, 1, $false, "24×19=", 2) | Out-Null
$d.Content.Find.Execute("46×49=", $true, $false, $false, $false, $false, $true, 1, $false, "42×20=", 2) | Out-Null
$d.Content.Find.Execute("84×77=", $true, $false, $false, $false, $false, $true, 1, $false, "35×94=", 2) | Out-Null
$d.Content.Find.Execute("89×33=", $true, $false, $false, $false, $false, $true, 1, $false, "41×92=", 2) | Out-Null
$d.Content.Find.Execute("60×39=", $true, $false, $false, $false, $false, $true, 1, $false, "35×55=", 2) | Out-Null
$d.Content.Find.Execute("17×56=", $true, $false, $false, $false, $false, $true, 1, $false, "72×35=", 2) | Out-Null
$d.Content.Find.Execute("62×90=", $true, $false, $false, $false, $false, $true, 1, $false, "74×14=", 2) | Out-Null
$d.Content.Find.Execute("95×28=", $true, $false, $false, $false, $false, $true, 1, $false, "48×50=", 2) | Out-Null
$d.Content.Find.Execute("36×37=", $true, $false, $false, $false, $false, $true, 1, $false, "98×47=", 2) | Out-Null
$d.Content.Find.Execute("35×98=", $true, $false, $false, $false, $false, $true, 1, $false, "77×64=", 2) | Out-Null
$d.Content.Find.Execute("27×76=", $true, $false, $false, $false, $false, $true, 1, $false, "34×21=", 2) | Out-Null
